$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @(0.6753301551942219, 1.667794583268128, 0.1575252929769615, 0.496779210170732, 2.997429241610044)
    3 = @(3.230985683306322, 1.667794583268128, 0.8054896365839992, 0.496779210170732, 6.201049113329182)
    4 = @(3.230985683306322, 1.667794583268128, 3.900430680208489, 0.496779210170732, 9.295990156953671)
    5 = @(3.230985683306322, 1.667794583268128, 0.8054896365839992, 8.660232485948974, 14.36450238910742)
    6 = @(3.230985683306322, 1.667794583268128, 0.8054896365839992, 0.496779210170732, 6.201049113329182)
    7 = @(3.230985683306322, 1.667794583268128, 0.8054896365839992, 0.496779210170732, 6.201049113329182)
    8 = @(3.230985683306322, 1.667794583268128, 0.1575252929769615, 0.496779210170732, 5.553084769722144)
    9 = @(3.230985683306322, 1.667794583268128, 0.8054896365839992, 0.496779210170732, 6.201049113329182)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
    $ws.Cells.Item($row, 5).Value = $vals[3]
    $ws.Cells.Item($row, 7).Value = $vals[4]
}
